# padel.xlsx update — record the new 1a vuelta match "Alba-Luis" vs "Teresa-Leticia"
# (6-0,6-1 / 0-6,1-6, Alba-Luis win) played on 2025-12-04, and roll the whole
# historial_partidos date column forward from 2025-12-01 to 2025-12-04.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) resultados — append the new match result as row 34
# ---------------------------------------------------------------------------
$wsResultados = $wb.Worksheets.Item("resultados")

$wsResultados.Range("A34").Value = "Mediocre medio"
$wsResultados.Range("B34").Value = "1ª vuelta"
$wsResultados.Range("C34").Value = "Alba-Luis"
$wsResultados.Range("D34").Value = "Teresa-Leticia"
$wsResultados.Range("E34").Value = "6-0,6-1"
$wsResultados.Range("F34").Value = "0-6,1-6"

# ---------------------------------------------------------------------------
# 2) clasificacion_auto — update the standings for the two pairs involved
#    row 15 = Alba-Luis, row 17 = Teresa-Leticia
# ---------------------------------------------------------------------------
$wsClasifAuto = $wb.Worksheets.Item("clasificacion_auto")

$wsClasifAuto.Range("D15").Value = 9
$wsClasifAuto.Range("E15").Value = 4
$wsClasifAuto.Range("F15").Value = 3
$wsClasifAuto.Range("G15").Value = 0
$wsClasifAuto.Range("H15").Value = 1
$wsClasifAuto.Range("I15").Value = 6
$wsClasifAuto.Range("J15").Value = 2
$wsClasifAuto.Range("K15").Value = 40
$wsClasifAuto.Range("L15").Value = 19

$wsClasifAuto.Range("D17").Value = 3
$wsClasifAuto.Range("E17").Value = 4
$wsClasifAuto.Range("F17").Value = 1
$wsClasifAuto.Range("G17").Value = 0
$wsClasifAuto.Range("H17").Value = 3
$wsClasifAuto.Range("I17").Value = 2
$wsClasifAuto.Range("J17").Value = 6
$wsClasifAuto.Range("K17").Value = 18
$wsClasifAuto.Range("L17").Value = 39

# ---------------------------------------------------------------------------
# 3) clasificacion — same two pairs, mirrored (no JG/JP columns here)
#    row 8 = Alba-Luis, row 10 = Teresa-Leticia
# ---------------------------------------------------------------------------
$wsClasif = $wb.Worksheets.Item("clasificacion")

$wsClasif.Range("D8").Value = 9
$wsClasif.Range("E8").Value = 4
$wsClasif.Range("F8").Value = 3
$wsClasif.Range("G8").Value = 0
$wsClasif.Range("H8").Value = 1
$wsClasif.Range("I8").Value = 6
$wsClasif.Range("J8").Value = 2

$wsClasif.Range("E10").Value = 4
$wsClasif.Range("H10").Value = 3
$wsClasif.Range("J10").Value = 6

# ---------------------------------------------------------------------------
# 4) historial_partidos — shift every existing match date from 2025-12-01
#    (serial 45992) to 2025-12-04 (serial 45995), then append the two new
#    per-pair rows (66 & 67) describing the new match.
# ---------------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item("historial_partidos")

for ($r = 2; $r -le 65; $r++) {
    $wsHist.Cells.Item($r, 1).Value = 45995
}

# new rows need the same date number-format as the rest of column A
$wsHist.Range("A65").Copy()
$wsHist.Range("A66:A67").PasteSpecial(-4122)
$wsHist.Application.CutCopyMode = $false

$wsHist.Range("A66").Value = 45995
$wsHist.Range("B66").Value = "mediocre medio"
$wsHist.Range("C66").Value = "1ª vuelta"
$wsHist.Range("D66").Value = "Alba-Luis"
$wsHist.Range("E66").Value = "Gana"
$wsHist.Range("F66").Value = 2
$wsHist.Range("G66").Value = 0
$wsHist.Range("H66").Value = 3
$wsHist.Range("I66").Value = 4
$wsHist.Range("J66").Value = 9
$wsHist.Range("K66").Value = 3
$wsHist.Range("L66").Value = 0
$wsHist.Range("M66").Value = 1

$wsHist.Range("A67").Value = 45995
$wsHist.Range("B67").Value = "mediocre medio"
$wsHist.Range("C67").Value = "1ª vuelta"
$wsHist.Range("D67").Value = "Teresa-Leticia"
$wsHist.Range("E67").Value = "Pierde"
$wsHist.Range("F67").Value = 0
$wsHist.Range("G67").Value = 2
$wsHist.Range("H67").Value = 0
$wsHist.Range("I67").Value = 4
$wsHist.Range("J67").Value = 3
$wsHist.Range("K67").Value = 1
$wsHist.Range("L67").Value = 0
$wsHist.Range("M67").Value = 3

# ---------------------------------------------------------------------------
# 5) restore resultados as the active sheet with F34 selected, matching the
#    workbook's saved cursor position after the edit.
# ---------------------------------------------------------------------------
$wsResultados.Activate()
$wsResultados.Range("F34").Select()
